$d = $word.ActiveDocument

# --- 1. Insert a new "Meta description" paragraph right after the H1 title ---
$titlePara = $d.Paragraphs.Item(1)
[void]$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:r/>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' + `
  '<w:r><w:t>: Experience the spooky feel of Halloween in Bat Stax, a 5-reel slot game with medium volatility. Play for free now!</w:t></w:r>' + `
  '</w:p>'
[void]$metaPara.Range.InsertXML($metaXml)

# --- 2. Remove the duplicated bold "Play Bat Stax for Free..." paragraph near
#        the end of the document (keep the original H1 title at index 1). ---
for ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Play Bat Stax for Free: Exciting Halloween-themed Slot Game`r") {
        [void]$para.Range.Delete()
    }
}

# --- 3. Replace the italic "Experience the spooky feel..." paragraph text
#        with the new feature-image prompt (keeps the run's <w:i/> formatting). ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Experience the spooky feel of Halloween in Bat Stax, a 5-reel slot game with medium volatility. Play for free now!`r") {
        $r = $d.Range($para.Range.Start, $para.Range.End)
        $r.Text = "Create a colorful and eye-catching feature image for Bat Stax, following the prompt below: Prompt: Design a cartoon-style image featuring a happy Maya warrior with glasses. The feature image should include the Bat Stax logo, and the Maya warrior should be holding a witch's hat and surrounded by bats. The background should be a foggy cemetery at night, with a glowing Jack O'Lantern nearby. Make sure to include the colors predominant in Bat Stax: dark blue, black, and purple. The overall image should reflect the spooky yet playful theme of Halloween."
    }
}

Write-Host "Edit complete."
